# ESS-BoB ToDo.docx -- final cleanup edits:
#  1. Title: split "BoB" into its own run, wrapped in proofErr spellStart/spellEnd
#  2. "Zvika" in the Earth-to-GND bullet: split into its own run w/ proofErr
#  3. "vias" in the (struck-through) polygon-pour bullet: split into its own
#     run w/ proofErr, keeping the <w:strike/> formatting on every run
#  4. Reorder the last two "Chassis GND" bullets: the CSDJ-connectors bullet
#     moves above the "Add 2-3 more GND test points" bullet, and the
#     _GoBack bookmark moves with the text it used to sit beside.
#
# w:proofErr elements aren't exposed as a COM object, and Find/Replace only
# ever touches plain text -- so each edit below is done by building the
# literal run-level OOXML for the paragraph(s) and dropping it in with
# Range.InsertXML, which replaces the range's contents in place while any
# untouched wrapping (like <w:pPr>) is preserved only if we re-supply it too.

$d = $word.ActiveDocument

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# Unicode chars that must be byte-exact (typing them literally risks silent
# normalization, e.g. OHM SIGN -> GREEK CAPITAL OMEGA) -- build via [char].
$rsquo = [char]0x2019   # RIGHT SINGLE QUOTATION MARK  (’)
$ohm   = [char]0x2126   # OHM SIGN                     (Ω)

# ---------------------------------------------------------------------
# 1) Title paragraph: "ESS-BoB Board Design ToDo"
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titleXml = $pkgOpen + `
  '<w:p><w:pPr><w:pStyle w:val="Title"/></w:pPr>' + `
  '<w:r><w:t>ESS-</w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>BoB</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> Board Design ToDo</w:t></w:r>' + `
  '</w:p>' + $pkgClose
$titlePara.Range.InsertXML($titleXml)

# ---------------------------------------------------------------------
# 2) "Earth to GND connection ... (consult with Zvika)"
# ---------------------------------------------------------------------
$earthPara = $d.Content.Find
$earthRange = $d.Content
$earthRange.Find.Execute("Earth to GND connection*Zvika)", $true, $false, $true, $false, $false, $true, 1, $false, "", 0) | Out-Null
$earthP = $earthRange.Paragraphs(1)
$earthXml = $pkgOpen + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
  '<w:r><w:t xml:space="preserve">+ </w:t></w:r>' + `
  "<w:r><w:t xml:space=`"preserve`">Earth to GND connection: decide what${rsquo}s best and implement (consult with </w:t></w:r>" + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>Zvika</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t>)</w:t></w:r>' + `
  '</w:p>' + $pkgClose
$earthP.Range.InsertXML($earthXml)

# ---------------------------------------------------------------------
# 3) "Add polygon pour for GND ... lot of vias between them ..." (struck through)
# ---------------------------------------------------------------------
$viasRange = $d.Content
$viasRange.Find.Execute("Add polygon pour*vias between them*70$*).", $true, $false, $true, $false, $false, $true, 1, $false, "", 0) | Out-Null
$viasP = $viasRange.Paragraphs(1)
$viasXml = $pkgOpen + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:strike/></w:rPr></w:pPr>' + `
  "<w:r><w:rPr><w:strike/></w:rPr><w:t xml:space=`"preserve`">Add polygon pour for GND on top side with a lot of </w:t></w:r>" + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:rPr><w:strike/></w:rPr><w:t>vias</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  "<w:r><w:rPr><w:strike/></w:rPr><w:t xml:space=`"preserve`"> between them. If it still doesn${rsquo}t look like good ground, switch to 4-layer (extra 70`$).</w:t></w:r>" + `
  '</w:p>' + $pkgClose
$viasP.Range.InsertXML($viasXml)

# ---------------------------------------------------------------------
# 4) Reorder: "Connect chassis GND to CSDJ connectors ..." bullet now comes
#    before "Add 2-3 more GND test points.", and the _GoBack bookmark now
#    sits with the "Add 2-3..." text instead.
# ---------------------------------------------------------------------
$csdjRange = $d.Content
$csdjRange.Find.Execute("Connect chassis GND to CSDJ*solder jumpers", $true, $false, $true, $false, $false, $true, 1, $false, "", 0) | Out-Null
$csdjP = $csdjRange.Paragraphs(1)
$csdjStart = $csdjP.Range.Start

$gndRange = $d.Content
$gndRange.Find.Execute("Add 2-3 more GND test points.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$gndP = $gndRange.Paragraphs(1)
$gndEnd = $gndP.Range.End

$spanRange = $d.Range($csdjStart, $gndEnd)
$reorderXml = $pkgOpen + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
  '<w:r><w:t xml:space="preserve">+ </w:t></w:r>' + `
  "<w:r><w:t>Connect chassis GND to CSDJ connectors using 0-${ohm} resistors or solder jumpers</w:t></w:r>" + `
  '</w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
  '<w:r><w:t xml:space="preserve">+ </w:t></w:r>' + `
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
  '<w:bookmarkEnd w:id="0"/>' + `
  '<w:r><w:t>Add 2-3 more GND test points.</w:t></w:r>' + `
  '</w:p>' + $pkgClose
$spanRange.InsertXML($reorderXml)
